# "Document Update, StartStatus Rename"
#
# 1) Rename the "StartStatus" sheet to "Status".
# 2) Update its selection (last-used cell moves from J219 -> G216) and
#    make it the active/selected sheet/tab (it becomes the workbook's
#    active tab, replacing "_Schedule").
# 3) The volatile NOW() formula on "_Schedule"!B2 recalculates automatically
#    as part of the normal save/recalc cycle, refreshing the "Document Update"
#    timestamp.

$wb = $excel.ActiveWorkbook

# --- Rename sheet: StartStatus -> Status -----------------------------------
$statusSheet = $wb.Worksheets.Item("StartStatus")
$statusSheet.Name = "Status"

# --- Update the sheet's remembered selection --------------------------------
$statusSheet.Range("G216").Select()

# --- Make it the active sheet/tab (was "_Schedule") -------------------------
$statusSheet.Activate()
